$wb = $excel.ActiveWorkbook

# Sheet3 ("100 Australian species" / VSEARCH row 13)
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Range("D13").Value = 0.2
$ws.Range("E13").Value = 0.01492537313432836
$ws.Range("F13").Value = 0.02777777777777778
$ws.Range("G13").Value = 0.05747126436781608
$ws.Range("H13").Value = 0.2929292929292929

# Sheet6 ("Lutjanidae" / VSEARCH row 13)
$ws = $wb.Worksheets.Item("Sheet6")
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0.2
$ws.Range("F13").Value = 0.3333333333333334
$ws.Range("G13").Value = 0.5555555555555556
$ws.Range("H13").Value = 0.2592592592592592

# Sheet9 ("Rottnest" / VSEARCH row 13)
$ws = $wb.Worksheets.Item("Sheet9")
$ws.Range("D13").Value = 0.7352941176470589
$ws.Range("E13").Value = 0.2358490566037736
$ws.Range("F13").Value = 0.3571428571428571
$ws.Range("G13").Value = 0.5165289256198347
$ws.Range("H13").Value = 0.2307692307692308
